$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 7041.6665
$ws.Range("I82").Value = 7041.6665
$ws.Range("K82").Value = 21124.9995
$ws.Range("M82").Value = -20718.9995
$ws.Range("H85").Value = 7041.6665
$ws.Range("I85").Value = 7041.6665
$ws.Range("K85").Value = 21124.9995
$ws.Range("M85").Value = -19720.9995
$ws.Range("H92").Value = 520.61536
$ws.Range("I92").Value = 499.5
$ws.Range("K92").Value = 499.5
$ws.Range("M92").Value = 748.5
$ws.Range("H106").Value = 2123.3
$ws.Range("I106").Value = 2906
$ws.Range("K106").Value = 2906
$ws.Range("M106").Value = -2275
$ws.Range("H111").Value = 693.4167
$ws.Range("I111").Value = 620.0909
$ws.Range("K111").Value = 1860.2727
$ws.Range("M111").Value = 1206.7273
$ws.Range("H137").Value = 317087.16
$ws.Range("I137").Value = 1795.9584
$ws.Range("K137").Value = 5387.8752
$ws.Range("M137").Value = -2837.8752

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 550
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726
$ws.Range("H32").Value = 16507.365
$ws.Range("I32").Value = 18115.916
$ws.Range("J32").Value = 14236.471
$ws.Range("K32").Value = 18115.916
$ws.Range("L32").Value = 14236.471
$ws.Range("M32").Value = -17828.916
$ws.Range("N32").Value = -14810.471
$ws.Range("H61").Value = 1583.1875
$ws.Range("I61").Value = 1345.0714
$ws.Range("K61").Value = 1345.0714
$ws.Range("M61").Value = -1133.0714
$ws.Range("H74").Value = 1995.4828
$ws.Range("I74").Value = 1553.1111
$ws.Range("J74").Value = 2719.3635
$ws.Range("K74").Value = 1553.1111
$ws.Range("L74").Value = 2719.3635
$ws.Range("M74").Value = -679.1111000000001
$ws.Range("N74").Value = -4467.363499999999
$ws.Range("H77").Value = 1995.4828
$ws.Range("I77").Value = 1553.1111
$ws.Range("J77").Value = 2719.3635
$ws.Range("K77").Value = 7765.5555
$ws.Range("L77").Value = 13596.8175
$ws.Range("M77").Value = -3397.5555
$ws.Range("N77").Value = -22332.8175
$ws.Range("H102").Value = 138490
$ws.Range("I102").Value = 201239.6
$ws.Range("J102").Value = 33907.332
$ws.Range("K102").Value = 201239.6
$ws.Range("L102").Value = 33907.332
$ws.Range("M102").Value = -199617.6
$ws.Range("N102").Value = -37151.332
$ws.Range("H107").Value = 50369.25
$ws.Range("J107").Value = 50369.25
$ws.Range("L107").Value = 50369.25
$ws.Range("N107").Value = -58049.25
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180
$ws.Range("H116").Value = 550
$ws.Range("J116").Value = 500
$ws.Range("L116").Value = 500
$ws.Range("N116").Value = -5088
$ws.Range("H122").Value = 2078.8518
$ws.Range("I122").Value = 1975.4783
$ws.Range("K122").Value = 5926.4349
$ws.Range("M122").Value = -3476.4349
$ws.Range("H132").Value = 1853
$ws.Range("J132").Value = 2357
$ws.Range("L132").Value = 7071
$ws.Range("N132").Value = -12131
$ws.Range("H136").Value = 1583.1875
$ws.Range("I136").Value = 1345.0714
$ws.Range("K136").Value = 4035.2142
$ws.Range("M136").Value = -1485.2142
$ws.Range("H138").Value = 70596.75
$ws.Range("J138").Value = 69332.336
$ws.Range("L138").Value = 69332.336
$ws.Range("N138").Value = -79612.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 550
$ws.Range("J3").Value = 500
$ws.Range("L3").Value = 500
$ws.Range("N3").Value = -728
$ws.Range("H86").Value = 5469.6787
$ws.Range("I86").Value = 4534.522
$ws.Range("J86").Value = 9771.4
$ws.Range("K86").Value = 4534.522
$ws.Range("L86").Value = 9771.4
$ws.Range("M86").Value = -3411.522
$ws.Range("N86").Value = -12017.4
$ws.Range("H89").Value = 5469.6787
$ws.Range("I89").Value = 4534.522
$ws.Range("J89").Value = 9771.4
$ws.Range("K89").Value = 22672.61
$ws.Range("L89").Value = 48857
$ws.Range("M89").Value = -17056.61
$ws.Range("N89").Value = -60089
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 1746
$ws.Range("H99").Value = 2483296.8
$ws.Range("I99").Value = 252494.75
$ws.Range("K99").Value = 252494.75
$ws.Range("M99").Value = -250996.75
$ws.Range("H107").Value = 2224.6667
$ws.Range("I107").Value = 1776.25
$ws.Range("J107").Value = 4018.3333
$ws.Range("K107").Value = 1776.25
$ws.Range("L107").Value = 4018.3333
$ws.Range("M107").Value = 143.75
$ws.Range("N107").Value = -7858.3333
$ws.Range("H134").Value = 3027.0466
$ws.Range("I134").Value = 2443
$ws.Range("K134").Value = 7329
$ws.Range("M134").Value = -4794
$ws.Range("H140").Value = 43500
$ws.Range("J140").Value = 43500
$ws.Range("L140").Value = 43500
$ws.Range("N140").Value = -53860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 136664.67
$ws.Range("J59").Value = 136664.67
$ws.Range("L59").Value = 136664.67
$ws.Range("N59").Value = -138954.67
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 55000
$ws.Range("J80").Value = 55000
$ws.Range("L80").Value = 55000
$ws.Range("N80").Value = -57246
$ws.Range("H83").Value = 55000
$ws.Range("J83").Value = 55000
$ws.Range("L83").Value = 165000
$ws.Range("N83").Value = -176232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 91831.17999999999
$ws.Range("J113").Value = 112010.336
$ws.Range("L113").Value = 336031.008
$ws.Range("N113").Value = -340371.008
$ws.Range("H122").Value = 1263762.8
$ws.Range("I122").Value = 295
$ws.Range("J122").Value = 1444258.1
$ws.Range("K122").Value = 2655
$ws.Range("L122").Value = 12998322.9
$ws.Range("M122").Value = -205
$ws.Range("N122").Value = -13003222.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1566.5454
$ws.Range("I102").Value = 1517.6666
$ws.Range("J102").Value = 1786.5
$ws.Range("K102").Value = 1517.6666
$ws.Range("L102").Value = 1786.5
$ws.Range("M102").Value = 104.3334
$ws.Range("N102").Value = -5030.5
$ws.Range("H122").Value = 585912.75
$ws.Range("I122").Value = 818453.2
$ws.Range("J122").Value = 4561.6665
$ws.Range("K122").Value = 2455359.6
$ws.Range("L122").Value = 13684.9995
$ws.Range("M122").Value = -2452909.6
$ws.Range("N122").Value = -18584.9995
$ws.Range("H138").Value = 137500
$ws.Range("J138").Value = 137500
$ws.Range("L138").Value = 137500
$ws.Range("N138").Value = -147780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8781.091
$ws.Range("J62").Value = 4430.6
$ws.Range("L62").Value = 4430.6
$ws.Range("N62").Value = -5678.6
$ws.Range("H65").Value = 8781.091
$ws.Range("J65").Value = 4430.6
$ws.Range("L65").Value = 22153
$ws.Range("N65").Value = -28393
$ws.Range("H107").Value = 14044.333
$ws.Range("I107").Value = 22739
$ws.Range("J107").Value = 3176
$ws.Range("K107").Value = 68217
$ws.Range("L107").Value = 9528
$ws.Range("M107").Value = -66297
$ws.Range("N107").Value = -13368
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 50000
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 50000
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -45100
$ws.Range("N123").ClearContents()
